# feat: add 2022-Q1 data
#
# 1. Clone the "2021-Q4" sheet (same column layout) into a new sheet placed
#    right before "总计", rename it "2022-Q1", trim it down to the 9 funds
#    that held the stock in 2022-Q1, and overwrite every cell with the
#    2022-Q1 figures.
# 2. Insert a new summary row at the top of "总计" for 2022-Q1 and
#    renumber the existing index column.
#
# NOTE: sheet object references captured in a variable can resolve to the
# wrong worksheet once a sheet is inserted/copied elsewhere in the
# workbook (the handle seems to be positional), so we always re-look the
# sheet up by its tab name with $wb.Worksheets.Item(...) right before each
# use instead of reusing a variable across a Copy()/Add()/Delete() call.

$wb = $excel.ActiveWorkbook

# --- 1. Build the "2022-Q1" sheet -----------------------------------------
$wb.Worksheets.Item("2021-Q4").Copy($wb.Worksheets.Item("总计"))
$wb.Worksheets.Item("2021-Q4 (2)").Name = "2022-Q1"

$newSheet = $wb.Worksheets.Item("2022-Q1")

# The template sheet has 19 data rows (rows 2-20); 2022-Q1 only needs 9.
$newSheet.Range("A11:A20").EntireRow.Delete()

$newSheet.Range("B1").Value2 = "基金代码"
$newSheet.Range("C1").Value2 = "基金名称"
$newSheet.Range("D1").Value2 = "基金规模"
$newSheet.Range("E1").Value2 = "股票总仓位"
$newSheet.Range("F1").Value2 = "仓位占比"
$newSheet.Range("G1").Value2 = "持有市值(亿元)"
$newSheet.Range("H1").Value2 = "仓位排名"

# index, 基金代码, 基金名称, 基金规模, 股票总仓位, 仓位占比, 持有市值(亿元), 仓位排名
$fundRows = @(
    @(0, "005368", "富国清洁能源产业灵活配置混合A",     "23.52", "88.60", "2.86", "0.6727", 10),
    @(1, "100029", "富国天成红利混合",                   "10.21", "76.00", "2.53", "0.2583", 6),
    @(2, "410001", "华富竞争力优选混合",                 "3.34",  "89.23", "6.28", "0.2098", 5),
    @(3, "000849", "汇丰晋信双核策略混合A",               "3.00",  "85.44", "3.86", "0.1158", 6),
    @(4, "014663", "富国创新发展两年定期开放混合A",       "2.62",  "37.11", "1.84", "0.0482", 6),
    @(5, "011127", "富国清洁能源产业灵活配置混合C",       "1.61",  "88.60", "2.86", "0.0460", 10),
    @(6, "007713", "华富科技动能混合",                   "0.56",  "86.98", "6.84", "0.0383", 5),
    @(7, "000850", "汇丰晋信双核策略混合C",               "0.39",  "85.44", "3.86", "0.0151", 6),
    @(8, "014664", "富国创新发展两年定期开放混合C",       "0.32",  "37.11", "1.84", "0.0059", 6)
)

$r = 2
foreach ($fund in $fundRows) {
    $newSheet.Range("A$r").Value2 = $fund[0]

    # Fund code, size, position%, weight%, and rank-weight are stored as
    # text in the source data (leading zeros / fixed decimals), so force
    # text entry and then drop the leftover number-format style.
    $newSheet.Range("B$r").NumberFormat = "@"
    $newSheet.Range("B$r").Value2 = $fund[1]
    $newSheet.Range("B$r").ClearFormats()

    $newSheet.Range("C$r").Value2 = $fund[2]

    $newSheet.Range("D$r").NumberFormat = "@"
    $newSheet.Range("D$r").Value2 = $fund[3]
    $newSheet.Range("D$r").ClearFormats()

    $newSheet.Range("E$r").NumberFormat = "@"
    $newSheet.Range("E$r").Value2 = $fund[4]
    $newSheet.Range("E$r").ClearFormats()

    $newSheet.Range("F$r").NumberFormat = "@"
    $newSheet.Range("F$r").Value2 = $fund[5]
    $newSheet.Range("F$r").ClearFormats()

    $newSheet.Range("G$r").NumberFormat = "@"
    $newSheet.Range("G$r").Value2 = $fund[6]
    $newSheet.Range("G$r").ClearFormats()

    $newSheet.Range("H$r").Value2 = $fund[7]

    $r = $r + 1
}

# --- 2. Update the "总计" (summary) sheet ---------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("A2:D2").ClearFormats()

$wb.Worksheets.Item("2021-Q4").Range("A2").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Range("A2").Value2 = 0
$totalSheet.Range("B2").Value2 = "2022-Q1"
$totalSheet.Range("C2").Value2 = 9
$totalSheet.Range("D2").Value2 = 1.41

$totalSheet.Range("A3").Value2 = 1
$totalSheet.Range("A4").Value2 = 2
